$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos price/volume refresh (GitHub Actions data pull)

$ws.Range("D2").Value = "27.726.16"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "1.905.26"
$ws.Range("E3").Value = "  +0.63%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9993"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5193"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.79%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3779"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.42%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07245"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.16"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9009"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.82%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07667"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "1.914.73"
$ws.Range("E13").Value = "  +0.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.440"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.39%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.00"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008714"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.000"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.16%  "
$ws.Range("D19").Value = "27.760.83"
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.50"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.44%  "
$ws.Range("E21").Value = "  +0.42%  "
$ws.Range("D22").Value = "2.110.01"
$ws.Range("E22").Value = "  -1.95%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.84"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.87%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.621"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.871"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.89%  "
$ws.Range("E27").Value = "  -0.34%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.159"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.49"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.855"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09030"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.95%  "
$ws.Range("E32").Value = "  -0.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.836"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.63%  "
$ws.Range("E34").Value = "  +0.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7805"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.43%  "
$ws.Range("E36").Value = "  +2.65%  "
$ws.Range("E37").Value = "  +1.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.071"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.79%  "
$ws.Range("E39").Value = "  -0.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5566"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.05283"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.714"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "114.77"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.77%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.525"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1517"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4808"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.50"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9998"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.614"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.90%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "66.71"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.92%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06000"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.08%  "
